# Insert a new weekly price row for "Ajo" (Chino, Primera) from Femacal de
# La Calera, dated 2021-09-22 (serial 44461), above the existing row 202.
# This pushes the previous rows 202-216 down to 203-217 (dimension grows
# from R216 to R217), matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 202, shifting rows 202:216 -> 203:217
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new data point
$ws.Cells.Item(202, 1).Value  = 3
$ws.Cells.Item(202, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(202, 3).Value  = "Coquimbo"
$ws.Cells.Item(202, 4).Value  = 44461
$ws.Cells.Item(202, 5).Value  = 5
$ws.Cells.Item(202, 6).Value  = 100112003
$ws.Cells.Item(202, 7).Value  = "Ajo"
$ws.Cells.Item(202, 8).Value  = "Chino"
$ws.Cells.Item(202, 9).Value  = "Primera"
$ws.Cells.Item(202, 10).Value = 65
$ws.Cells.Item(202, 11).Value = 16000
$ws.Cells.Item(202, 12).Value = 16500
$ws.Cells.Item(202, 13).Value = 16231
$ws.Cells.Item(202, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(202, 15).Value = "China"
$ws.Cells.Item(202, 16).Value = 1623
$ws.Cells.Item(202, 17).Value = 10
$ws.Cells.Item(202, 18).Value = "Hortaliza"
